$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.425.14'
$ws.Range('E2').Value = '  +1.70%  '
$ws.Range('D3').Value = '1.887.87'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.17'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.691'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.19'
$ws.Range('E8').Value = '  +2.91%  '
$ws.Range('E9').Value = '  +2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.61'
$ws.Range('E10').Value = '  +6.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0745'
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0987'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.89'
$ws.Range('E13').Value = '  +7.60%  '
$ws.Range('D14').Value = '2.162.75'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.771'
$ws.Range('E15').Value = '  +7.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.01'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '1.914.48'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('D18').Value = '35.443.85'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.47'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = '0.0₃0826'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '245.30'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.84'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.16'
$ws.Range('E23').Value = '  +4.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.65'
$ws.Range('E24').Value = '  +8.16%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.38'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.61'
$ws.Range('E28').Value = '  +2.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.30'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0596'
$ws.Range('E31').Value = '  +3.27%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('E33').Value = '  +24.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.21'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('E36').Value = '  -13.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.856'
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.94'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0724'
$ws.Range('E39').Value = '  +10.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0222'
$ws.Range('E40').Value = '  +5.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.30'
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.25'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.70'
$ws.Range('E44').Value = '  +12.44%  '
$ws.Range('D45').Value = '1.331.61'
$ws.Range('E45').Value = '  +3.70%  '
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0808'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.30'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').Value = '2.064.46'
$ws.Range('E51').Value = '  +0.30%  '
